$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new value would otherwise be
# auto-interpreted as a number by Excel, so they remain text like the source data.
$textCells = @('D4', 'D5', 'D6', 'D8', 'D10', 'D11', 'D16', 'D17', 'D18', 'D19', 'D23', 'D25', 'D26', 'D29', 'D32', 'D34', 'D35', 'D39', 'D40', 'D41', 'D42', 'D45', 'D48', 'D49', 'D50')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin data
$ws.Range('D2').Value = '28.434.10'
$ws.Range('E2').Value = '  -0.08%  '
$ws.Range('D3').Value = '1.549.08'
$ws.Range('E3').Value = '  -2.04%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.32%  '
$ws.Range('D5').Value = '210.60'
$ws.Range('E5').Value = '  -1.35%  '
$ws.Range('D6').Value = '0.481'
$ws.Range('E6').Value = '  -2.14%  '
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('D8').Value = '23.98'
$ws.Range('E8').Value = '  -0.22%  '
$ws.Range('E9').Value = '  -2.11%  '
$ws.Range('D10').Value = '0.0582'
$ws.Range('E10').Value = '  -1.63%  '
$ws.Range('D11').Value = '0.0889'
$ws.Range('E11').Value = '  -0.44%  '
$ws.Range('D12').Value = '1.771.58'
$ws.Range('E12').Value = '  -2.03%  '
$ws.Range('D13').Value = '1.551.97'
$ws.Range('E13').Value = '  -1.90%  '
$ws.Range('D14').Value = '28.420.06'
$ws.Range('E14').Value = '  -0.22%  '
$ws.Range('E15').Value = '  -2.04%  '
$ws.Range('D16').Value = '0.508'
$ws.Range('E16').Value = '  -2.21%  '
$ws.Range('D17').Value = '60.88'
$ws.Range('E17').Value = '  -2.11%  '
$ws.Range('D18').Value = '228.97'
$ws.Range('E18').Value = '  -0.71%  '
$ws.Range('D19').Value = '7.31'
$ws.Range('E19').Value = '  -1.54%  '
$ws.Range('E20').Value = '  -2.95%  '
$ws.Range('E21').Value = '  -0.27%  '
$ws.Range('E22').Value = '  -0.35%  '
$ws.Range('D23').Value = '8.90'
$ws.Range('E23').Value = '  -2.69%  '
$ws.Range('E24').Value = '  -2.37%  '
$ws.Range('D25').Value = '151.24'
$ws.Range('E25').Value = '  -0.34%  '
$ws.Range('D26').Value = '14.76'
$ws.Range('E26').Value = '  -1.85%  '
$ws.Range('E27').Value = '  -1.55%  '
$ws.Range('E28').Value = '  -0.26%  '
$ws.Range('D29').Value = '6.21'
$ws.Range('E29').Value = '  -3.73%  '
$ws.Range('E30').Value = '  -3.50%  '
$ws.Range('E31').Value = '  -4.69%  '
$ws.Range('D32').Value = '3.15'
$ws.Range('E32').Value = '  -1.82%  '
$ws.Range('D33').Value = '1.383.45'
$ws.Range('E33').Value = '  -1.09%  '
$ws.Range('D34').Value = '2.99'
$ws.Range('E34').Value = '  -3.71%  '
$ws.Range('D35').Value = '1.06'
$ws.Range('E35').Value = '  -2.49%  '
$ws.Range('E36').Value = '  -3.27%  '
$ws.Range('E37').Value = '  -2.55%  '
$ws.Range('E38').Value = '  -3.20%  '
$ws.Range('D39').Value = '0.0161'
$ws.Range('E39').Value = '  -2.68%  '
$ws.Range('D40').Value = '1.92'
$ws.Range('E40').Value = '  +1.32%  '
$ws.Range('B41').Value = 'ImmutableX'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D41').Value = '0.509'
$ws.Range('E41').Value = '  -2.46%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  -0.27%  '
$ws.Range('E43').Value = '  -2.59%  '
$ws.Range('E44').Value = '  -1.63%  '
$ws.Range('D45').Value = '5.31'
$ws.Range('E45').Value = '  -2.74%  '
$ws.Range('E46').Value = '  -2.51%  '
$ws.Range('D47').Value = '1.684.43'
$ws.Range('E47').Value = '  -2.03%  '
$ws.Range('D48').Value = '0.877'
$ws.Range('E48').Value = '  -8.68%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').Value = '85.66'
$ws.Range('E49').Value = '  -1.22%  '
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D50').Value = '43.30'
$ws.Range('E50').Value = '  +8.53%  '
$ws.Range('D51').Value = '0.0₆0101'
$ws.Range('E51').Value = '  -1.61%  '
